# Update localization Status for two files (b155ba78-... and e97209b6-...)
# from "Ready for handoff" to "In Translation" across all three sheets:
#   - Overview (per-language status columns E = zh-cn, F = de-de)
#   - zh-cn    (Status column C)
#   - de-de    (Status column C)

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus   # zh-cn status for b155ba78-008b-412b-b4bd-ed7c1c13672c.md
$wsOverview.Range("F3").Value = $newStatus   # de-de status for b155ba78-008b-412b-b4bd-ed7c1c13672c.md
$wsOverview.Range("E4").Value = $newStatus   # zh-cn status for e97209b6-bb91-485f-b088-944ad9623264.md
$wsOverview.Range("F4").Value = $newStatus   # de-de status for e97209b6-bb91-485f-b088-944ad9623264.md

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus       # b155ba78-008b-412b-b4bd-ed7c1c13672c.md
$wsZhCn.Range("C4").Value = $newStatus       # e97209b6-bb91-485f-b088-944ad9623264.md

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus       # b155ba78-008b-412b-b4bd-ed7c1c13672c.md
$wsDeDe.Range("C4").Value = $newStatus       # e97209b6-bb91-485f-b088-944ad9623264.md
